$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.299610048532486
$ws.Range("B1").Value = 2.469788551330566
$ws.Range("C1").Value = 4.11912202835083
$ws.Range("D1").Value = 1.453192591667175
$ws.Range("E1").Value = 0.6663867235183716
